# Weekly update: insert the newest "Arveja Verde" price record as a new
# row 18, pushing the existing rows 18-39 down to 19-40 (no other data
# changes — everything else just shifts down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18; Excel automatically shifts all
# rows currently at 18..39 down to 19..40 and extends the used range.
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with this week's record.
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C18").Value = "Ñuble"
$ws.Range("D18").Value = 44546
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 100112022
$ws.Range("G18").Value = "Arveja Verde"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 13000
$ws.Range("L18").Value = 14000
$ws.Range("M18").Value = 13500
$ws.Range("N18").Value = "$/saco 25 kilos"
$ws.Range("O18").Value = "Provincia de Diguillín"
$ws.Range("P18").Value = 540
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = "Hortaliza"
